$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 - Art_023 - puddingmv
$ws.Range("B24").Value = "https://puddingmv.tumblr.com/post/643034001960009728/she-knelt-to-the-ground-and-pulled-out-a-ring-and"
$ws.Range("D24").Value = "puddingmv"
$ws.Range("C24").Value = "?"
$ws.Range("E24").Value = "Tumblr"

# Row 25 - Art_024 - scavenger-rey
$ws.Range("B25").Value = "https://scavenger-rey.tumblr.com/post/643016761353125888/stargazing-date"
$ws.Range("D25").Value = "scavenger-rey"
$ws.Range("C25").Value = "?"
$ws.Range("E25").Value = "Tumblr"

# Row 26 - Art_025 - _camille_pelletier
$ws.Range("B26").Value = "https://www.instagram.com/p/B_29Awenjv3/"
$ws.Range("D26").Value = "_camille_pelletier"
$ws.Range("C26").Value = "Camille "
$ws.Range("E26").Value = "IG"

# Row 28 - Art_027 - brunamz
$ws.Range("B28").Value = "https://brunamz.tumblr.com/post/642424831199051776"
$ws.Range("D28").Value = "brunamz"
$ws.Range("C28").Value = "?"
$ws.Range("E28").Value = "Tumblr"

# update selection
$ws.Range("F17").Select()
